# Updated cryptos list on Mon Aug  7 09:41:57 UTC 2023 with GitHub Actions
#
# All Price (D) / Volume (E) cells are stored as plain text in this sheet
# (prices use '.' as both thousands- and decimal-separator, so a value like
# "29.056.75" is not a legal number anyway; other prices like "0.9988" ARE
# legal numbers and would silently be re-typed by Excel's auto-detection,
# losing trailing zeros / exact formatting). Using Value2 keeps whitespace-
# padded percentage strings as text automatically; for Price cells that
# look like plain numbers we add a leading apostrophe (Excel's standard
# "force text" marker) so the literal digits/zeros survive untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value2 = "29.056.75"
$ws.Range("E2").Value2 = "  +0.11%  "

# Row 3 - Ethereum
$ws.Range("D3").Value2 = "1.833.68"
$ws.Range("E3").Value2 = "  +0.32%  "

# Row 4 - TetherUSD
$ws.Range("D4").Value2 = "'0.9988"
$ws.Range("E4").Value2 = "  -0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value2 = "'242.85"
$ws.Range("E5").Value2 = "  -0.63%  "

# Row 6 - XRP
$ws.Range("D6").Value2 = "'0.6197"
$ws.Range("E6").Value2 = "  -1.99%  "

# Row 7 - USDC
$ws.Range("D7").Value2 = "'1.000"
$ws.Range("E7").Value2 = "  +0.10%  "

# Row 8 - Dogecoin
$ws.Range("D8").Value2 = "'0.07473"
$ws.Range("E8").Value2 = "  -0.74%  "

# Row 9 - Cardano
$ws.Range("D9").Value2 = "'0.2929"
$ws.Range("E9").Value2 = "  -0.41%  "

# Row 10 - Solana
$ws.Range("D10").Value2 = "'23.11"
$ws.Range("E10").Value2 = "  +0.17%  "

# Row 11 - TRON
$ws.Range("D11").Value2 = "'0.07676"
$ws.Range("E11").Value2 = "  -0.30%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value2 = "1.827.32"
$ws.Range("E12").Value2 = "  -0.02%  "

# Row 13 - Polkadot
$ws.Range("D13").Value2 = "'5.006"
$ws.Range("E13").Value2 = "  +0.18%  "

# Row 14 - Polygon
$ws.Range("D14").Value2 = "'0.6733"
$ws.Range("E14").Value2 = "  +0.62%  "

# Row 15 - Litecoin
$ws.Range("D15").Value2 = "'82.89"
$ws.Range("E15").Value2 = "  -0.32%  "

# Row 16 - ShibaInu
$ws.Range("D16").Value2 = "'0.000009164"
$ws.Range("E16").Value2 = "  -4.44%  "

# Row 17 - Uniswap
$ws.Range("D17").Value2 = "'5.915"
$ws.Range("E17").Value2 = "  -2.57%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value2 = "29.065.63"
$ws.Range("E18").Value2 = "  +0.05%  "

# Row 19 - WrappedliquidstakedEther2.0
$ws.Range("D19").Value2 = "2.080.04"
$ws.Range("E19").Value2 = "  +0.39%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value2 = "'239.77"
$ws.Range("E20").Value2 = "  +5.75%  "

# Row 21 - Avalanche
$ws.Range("D21").Value2 = "'12.69"
$ws.Range("E21").Value2 = "  +0.77%  "

# Row 22 - Dai (only Volume changes)
$ws.Range("E22").Value2 = "  +0.19%  "

# Row 23 - Chainlink
$ws.Range("D23").Value2 = "'7.206"
$ws.Range("E23").Value2 = "  +0.88%  "

# Row 24 - BinanceUSD
$ws.Range("D24").Value2 = "'1.0000"
$ws.Range("E24").Value2 = "  +0.00%  "

# Row 25 - Monero
$ws.Range("D25").Value2 = "'159.29"
$ws.Range("E25").Value2 = "  -0.65%  "

# Row 26 - Stellar
$ws.Range("D26").Value2 = "'0.1428"
$ws.Range("E26").Value2 = "  +0.12%  "

# Row 27 - Cosmos
$ws.Range("D27").Value2 = "'8.503"
$ws.Range("E27").Value2 = "  -0.12%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value2 = "'17.89"
$ws.Range("E28").Value2 = "  -0.21%  "

# Row 29 - PancakeSwap
$ws.Range("D29").Value2 = "'1.499"
$ws.Range("E29").Value2 = "  -0.39%  "

# Row 30 - Filecoin
$ws.Range("D30").Value2 = "'4.151"
$ws.Range("E30").Value2 = "  +0.07%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value2 = "'4.121"
$ws.Range("E31").Value2 = "  +1.39%  "

# Row 32 - Hedera
$ws.Range("D32").Value2 = "'0.05567"
$ws.Range("E32").Value2 = "  +1.66%  "

# Row 33 - Toncoin
$ws.Range("D33").Value2 = "'1.207"
$ws.Range("E33").Value2 = "  +0.50%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").Value2 = "'1.841"
$ws.Range("E34").Value2 = "  -0.80%  "

# Row 35 - ImmutableX
$ws.Range("D35").Value2 = "'0.7377"
$ws.Range("E35").Value2 = "  -0.99%  "

# Row 36 - ARBITRUM
$ws.Range("D36").Value2 = "'1.140"
$ws.Range("E36").Value2 = "  +0.31%  "

# Row 37 - HuobiToken
$ws.Range("D37").Value2 = "'2.658"
$ws.Range("E37").Value2 = "  +0.08%  "

# Row 38 - MXToken
$ws.Range("D38").Value2 = "'2.773"
$ws.Range("E38").Value2 = "  +0.66%  "

# Row 39 - VeChain
$ws.Range("D39").Value2 = "'0.01781"
$ws.Range("E39").Value2 = "  -0.15%  "

# Row 40 - Maker
$ws.Range("D40").Value2 = "1.211.02"
$ws.Range("E40").Value2 = "  -2.65%  "

# Row 41 - FraxShare
$ws.Range("D41").Value2 = "'6.476"
$ws.Range("E41").Value2 = "  -2.22%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").Value2 = "'0.8965"
$ws.Range("E42").Value2 = "  -0.64%  "

# Row 43 - PaxDollar
$ws.Range("D43").Value2 = "'0.9997"
$ws.Range("E43").Value2 = "  +0.06%  "

# Row 44 - Quant
$ws.Range("D44").Value2 = "'101.77"
$ws.Range("E44").Value2 = "  +0.54%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value2 = "1.977.64"
$ws.Range("E45").Value2 = "  +0.10%  "

# Row 46 - Aave
$ws.Range("D46").Value2 = "'65.69"
$ws.Range("E46").Value2 = "  +1.06%  "

# Rows 47 & 48 swap places: Mantle moves up to rank 47, BabyDogeCoin drops to rank 48
$ws.Range("B47").Value2 = "Mantle"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value2 = "'0.5083"
$ws.Range("E47").Value2 = "  -0.35%  "

$ws.Range("B48").Value2 = "BabyDogeCoin"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value2 = "'0.00000000119"
$ws.Range("E48").Value2 = "  -3.92%  "

# Row 49 - TheSandbox
$ws.Range("D49").Value2 = "'0.4073"
$ws.Range("E49").Value2 = "  +0.16%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value2 = "'9.136"
$ws.Range("E50").Value2 = "  +1.44%  "

# Row 51 - Cronos
$ws.Range("D51").Value2 = "'0.05819"
$ws.Range("E51").Value2 = "  +0.58%  "
